$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.782.45"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "2.281.68"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'251.76"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'0.639"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("D7").Value = "'74.68"
$ws.Range("E7").Value = "  +5.98%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.646"
$ws.Range("E9").Value = "  -2.93%  "
$ws.Range("D10").Value = "'39.81"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").Value = "'0.0974"
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").Value = "'7.50"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +1.58%  "
$ws.Range("D14").Value = "2.620.99"
$ws.Range("D15").Value = "'15.07"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'0.873"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "2.280.89"
$ws.Range("E17").Value = "  +1.74%  "
$ws.Range("D18").Value = "42.653.23"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'72.48"
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "'235.33"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +5.20%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'11.35"
$ws.Range("E26").Value = "  -1.66%  "
$ws.Range("D27").Value = "'2.40"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  +2.81%  "
$ws.Range("D29").Value = "'167.85"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'21.05"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("D31").Value = "'0.0887"
$ws.Range("E31").Value = "  +12.05%  "
$ws.Range("D32").Value = "'6.30"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").Value = "'0.127"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'31.65"
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "'4.59"
$ws.Range("E36").Value = "  +4.03%  "
$ws.Range("D37").Value = "'4.78"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "'0.0306"
$ws.Range("E38").Value = "  -4.15%  "
$ws.Range("D39").Value = "'13.70"
$ws.Range("E39").Value = "  +9.42%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'5.89"
$ws.Range("E41").Value = "  +1.59%  "
$ws.Range("D42").Value = "'0.211"
$ws.Range("E42").Value = "  +4.79%  "
$ws.Range("D43").Value = "'9.03"
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "'61.31"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("D45").Value = "'4.79"
$ws.Range("D46").Value = "'105.02"
$ws.Range("E46").Value = "  +10.75%  "
$ws.Range("E47").Value = "  -1.10%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'4.22"
$ws.Range("E51").Value = "  -0.99%  "
